# Update "想去人数" (interested count) values in column F for rows 2,3,5,6,7,8
# on both the "展览" and "全部类型" worksheets, per the source data refresh.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2" = 2344
    "F3" = 1850
    "F5" = 1136
    "F6" = 1093
    "F7" = 46
    "F8" = 5946
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
